# heuristic fixed, direction arrows added
# Update the demand/cost heuristic values for the three "locker" rows whose
# G column had been left at 1 (now corrected to 10).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G14").Value = 10
$ws.Range("G18").Value = 10
$ws.Range("G21").Value = 10

# Reflect the final active selection used while reviewing the fix.
$ws.Range("G6").Select()
